$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# 1. Duplicate "Datos Personales" sheet to create the new "Consulta Inicial" sheet
$ws1.Copy([System.Reflection.Missing]::Value, $ws1) | Out-Null
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Consulta Inicial"

# 2. Make room for the two new rows (Covid / Fecha de Covid) right before the
#    last row ("Otros"), which slides from row 16 down to row 18.
$ws2.Rows.Item(16).Insert()
$ws2.Rows.Item(16).Insert()

# 3. Update the text labels for the new sheet (values are typed in an order
#    that reproduces the shared-string table layout of the target workbook).
$ws2.Range("B2").Value  = "Consulta incial del paciente"
$ws2.Range("B8").Value  = "Motivo"
$ws2.Range("B7").Value  = "Fecha de Consulta"
$ws2.Range("B9").Value  = "Actividad Física"
$ws2.Range("B10").Value = "Antigüedad"
$ws2.Range("B11").Value = "Localización"
$ws2.Range("B12").Value = "Intensidad"
$ws2.Range("B13").Value = "Característica"
$ws2.Range("B14").Value = "Irradiación"
$ws2.Range("B15").Value = "Atenua"
$ws2.Range("B16").Value = "Covid"
$ws2.Range("B17").Value = "Fecha de Covid"

# 4. Merge the label/value cells of the two newly inserted rows like the
#    other data rows on the sheet.
$ws2.Range("B16:C16").Merge()
$ws2.Range("D16:F16").Merge()
$ws2.Range("B17:C17").Merge()
$ws2.Range("D17:F17").Merge()

# 5. Give the new sheet its own accent colour (Accent4 theme, lightened) so
#    it is visually distinct from "Datos Personales".
$used = $ws2.Range("B2:F18")
$used.Interior.ThemeColor = 8

# 6. Restore a left border on the value cells of the two new rows (matching
#    the look of the other single-row fields).
$ws2.Range("D16:F16").Borders.Item(7).LineStyle = 1
$ws2.Range("D16:F16").Borders.Item(7).Weight = 2
$ws2.Range("D17:F17").Borders.Item(7).LineStyle = 1
$ws2.Range("D17:F17").Borders.Item(7).Weight = 2
$ws2.Range("D16:F17").HorizontalAlignment = -4108

# 7. The new sheet becomes the active tab, with the first data-entry cell
#    selected, and "Datos Personales" is no longer the tab shown on open.
$ws2.Select()
$ws2.Range("B10:C10").Select()

$wb.Worksheets.Item(1).Range("B7:F16").Select()
